$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.846571
$ws.Range("H2").Value = 8.539712999999999
$ws.Range("I2").Value = 0.01041928469143245
$ws.Range("J2").Value = 0.01041928469143244
$ws.Range("M2").Value = 55.848606
$ws.Range("N2").Value = 167.545818
$ws.Range("O2").Value = 0.2323375192077237
$ws.Range("P2").Value = 0.2323375192077236
$ws.Range("Q2").Value = 158.977022230026
$ws.Range("R2").Value = 1430.793200070234
$ws.Range("S2").Value = 0.002420790757126427
$ws.Range("T2").Value = 0.002420790757126426
$ws.Range("G3").Value = 2.846571
$ws.Range("H3").Value = 8.539712999999999
$ws.Range("I3").Value = 0.01041928469143245
$ws.Range("J3").Value = 0.01041928469143244
$ws.Range("O3").Value = 0.3515710112922583
$ws.Range("P3").Value = 0.3515710112922583
$ws.Range("Q3").Value = 240.5625775304569
$ws.Range("R3").Value = 2165.063197774113
$ws.Range("S3").Value = 0.00366311845590885
$ws.Range("T3").Value = 0.003663118455908849
$ws.Range("G4").Value = 2.846571
$ws.Range("H4").Value = 8.539712999999999
$ws.Range("I4").Value = 0.01041928469143245
$ws.Range("J4").Value = 0.01041928469143244
$ws.Range("M4").Value = 33.195992
$ws.Range("N4").Value = 99.58797600000001
$ws.Range("O4").Value = 0.1380996766314891
$ws.Range("P4").Value = 0.1380996766314891
$ws.Range("Q4").Value = 94.49474814343199
$ws.Range("R4").Value = 850.452733290888
$ws.Range("S4").Value = 0.001438899846618246
$ws.Range("T4").Value = 0.001438899846618245
$ws.Range("G5").Value = 2.846571
$ws.Range("H5").Value = 8.539712999999999
$ws.Range("I5").Value = 0.01041928469143245
$ws.Range("J5").Value = 0.01041928469143244
$ws.Range("M5").Value = 66.82284533333335
$ws.Range("N5").Value = 200.468536
$ws.Range("O5").Value = 0.277991792868529
$ws.Range("P5").Value = 0.2779917928685289
$ws.Range("Q5").Value = 190.215973663352
$ws.Range("R5").Value = 1711.943762970168
$ws.Range("S5").Value = 0.002896475631778923
$ws.Range("T5").Value = 0.002896475631778922
$ws.Range("I6").Value = 0.07835537840126532
$ws.Range("J6").Value = 0.0783553784012653
$ws.Range("M6").Value = 55.848606
$ws.Range("N6").Value = 167.545818
$ws.Range("O6").Value = 0.2323375192077237
$ws.Range("P6").Value = 0.2323375192077236
$ws.Range("Q6").Value = 1195.543178139948
$ws.Range("R6").Value = 10759.88860325953
$ws.Range("S6").Value = 0.01820489423433243
$ws.Range("T6").Value = 0.01820489423433243
$ws.Range("I7").Value = 0.07835537840126532
$ws.Range("J7").Value = 0.0783553784012653
$ws.Range("O7").Value = 0.3515710112922583
$ws.Range("P7").Value = 0.3515710112922583
$ws.Range("S7").Value = 0.02754747962472042
$ws.Range("T7").Value = 0.02754747962472041
$ws.Range("I8").Value = 0.07835537840126532
$ws.Range("J8").Value = 0.0783553784012653
$ws.Range("M8").Value = 33.195992
$ws.Range("N8").Value = 99.58797600000001
$ws.Range("O8").Value = 0.1380996766314891
$ws.Range("P8").Value = 0.1380996766314891
$ws.Range("Q8").Value = 710.6218869131361
$ws.Range("R8").Value = 6395.596982218225
$ws.Range("S8").Value = 0.01082085241955271
$ws.Range("T8").Value = 0.0108208524195527
$ws.Range("I9").Value = 0.07835537840126532
$ws.Range("J9").Value = 0.0783553784012653
$ws.Range("M9").Value = 66.82284533333335
$ws.Range("N9").Value = 200.468536
$ws.Range("O9").Value = 0.277991792868529
$ws.Range("P9").Value = 0.2779917928685289
$ws.Range("Q9").Value = 1430.46716120663
$ws.Range("R9").Value = 12874.20445085967
$ws.Range("S9").Value = 0.02178215212265976
$ws.Range("T9").Value = 0.02178215212265975
$ws.Range("G10").Value = 7.347547
$ws.Range("H10").Value = 22.042641
$ws.Range("I10").Value = 0.02689417688042223
$ws.Range("J10").Value = 0.02689417688042223
$ws.Range("M10").Value = 55.848606
$ws.Range("N10").Value = 167.545818
$ws.Range("O10").Value = 0.2323375192077237
$ws.Range("P10").Value = 0.2323375192077236
$ws.Range("Q10").Value = 410.350257469482
$ws.Range("R10").Value = 3693.152317225338
$ws.Range("S10").Value = 0.006248526337531019
$ws.Range("T10").Value = 0.006248526337531017
$ws.Range("G11").Value = 7.347547
$ws.Range("H11").Value = 22.042641
$ws.Range("I11").Value = 0.02689417688042223
$ws.Range("J11").Value = 0.02689417688042223
$ws.Range("O11").Value = 0.3515710112922583
$ws.Range("P11").Value = 0.3515710112922583
$ws.Range("Q11").Value = 620.9382604003822
$ws.Range("R11").Value = 5588.444343603441
$ws.Range("S11").Value = 0.009455212963722917
$ws.Range("T11").Value = 0.009455212963722915
$ws.Range("G12").Value = 7.347547
$ws.Range("H12").Value = 22.042641
$ws.Range("I12").Value = 0.02689417688042223
$ws.Range("J12").Value = 0.02689417688042223
$ws.Range("M12").Value = 33.195992
$ws.Range("N12").Value = 99.58797600000001
$ws.Range("O12").Value = 0.1380996766314891
$ws.Range("P12").Value = 0.1380996766314891
$ws.Range("Q12").Value = 243.909111431624
$ws.Range("R12").Value = 2195.182002884616
$ws.Range("S12").Value = 0.003714077130456381
$ws.Range("T12").Value = 0.00371407713045638
$ws.Range("G13").Value = 7.347547
$ws.Range("H13").Value = 22.042641
$ws.Range("I13").Value = 0.02689417688042223
$ws.Range("J13").Value = 0.02689417688042223
$ws.Range("M13").Value = 66.82284533333335
$ws.Range("N13").Value = 200.468536
$ws.Range("O13").Value = 0.277991792868529
$ws.Range("P13").Value = 0.2779917928685289
$ws.Range("Q13").Value = 490.9839967603974
$ws.Range("R13").Value = 4418.855970843577
$ws.Range("S13").Value = 0.007476360448711919
$ws.Range("T13").Value = 0.007476360448711916
$ws.Range("G14").Value = 241.601176
$ws.Range("H14").Value = 724.8035279999999
$ws.Range("I14").Value = 0.8843311600268801
$ws.Range("J14").Value = 0.8843311600268799
$ws.Range("M14").Value = 55.848606
$ws.Range("N14").Value = 167.545818
$ws.Range("O14").Value = 0.2323375192077237
$ws.Range("P14").Value = 0.2323375192077236
$ws.Range("Q14").Value = 13493.08888756065
$ws.Range("R14").Value = 121437.7999880459
$ws.Range("S14").Value = 0.2054633078787338
$ws.Range("T14").Value = 0.2054633078787337
$ws.Range("G15").Value = 241.601176
$ws.Range("H15").Value = 724.8035279999999
$ws.Range("I15").Value = 0.8843311600268801
$ws.Range("J15").Value = 0.8843311600268799
$ws.Range("O15").Value = 0.3515710112922583
$ws.Range("P15").Value = 0.3515710112922583
$ws.Range("Q15").Value = 20417.61882382332
$ws.Range("R15").Value = 183758.5694144099
$ws.Range("S15").Value = 0.3109052002479061
$ws.Range("T15").Value = 0.3109052002479061
$ws.Range("G16").Value = 241.601176
$ws.Range("H16").Value = 724.8035279999999
$ws.Range("I16").Value = 0.8843311600268801
$ws.Range("J16").Value = 0.8843311600268799
$ws.Range("M16").Value = 33.195992
$ws.Range("N16").Value = 99.58797600000001
$ws.Range("O16").Value = 0.1380996766314891
$ws.Range("P16").Value = 0.1380996766314891
$ws.Range("Q16").Value = 8020.190705686592
$ws.Range("R16").Value = 72181.71635117933
$ws.Range("S16").Value = 0.1221258472348618
$ws.Range("T16").Value = 0.1221258472348618
$ws.Range("G17").Value = 241.601176
$ws.Range("H17").Value = 724.8035279999999
$ws.Range("I17").Value = 0.8843311600268801
$ws.Range("J17").Value = 0.8843311600268799
$ws.Range("M17").Value = 66.82284533333335
$ws.Range("N17").Value = 200.468536
$ws.Range("O17").Value = 0.277991792868529
$ws.Range("P17").Value = 0.2779917928685289
$ws.Range("Q17").Value = 16144.47801619945
$ws.Range("R17").Value = 145300.302145795
$ws.Range("S17").Value = 0.2458368046653784
$ws.Range("T17").Value = 0.2458368046653783
